# Split the single run "do Cục CSQLHC về TTXH cấp" into three runs:
#   "do "  |  "Cục CSQLHC về TTXH"  |  " cấp"
# while keeping identical run formatting (color 000000) on all three.

$d = $word.ActiveDocument

$target = "do Cục CSQLHC về TTXH cấp"
$part1  = "do "
$part2  = "Cục CSQLHC về TTXH"
$part3  = " cấp"

$rng = $d.Content
$found = $rng.Find.Execute($target, $true, $false, $false, $false, $false, `
                            $true, 1, $false, "", 0)

if ($found) {
    $start = $rng.Start

    # Boundaries between the three desired text chunks.
    $b1 = $start + $part1.Length
    $b2 = $b1 + $part2.Length

    # Force Word to split the run at these two boundaries by toggling a
    # character property on the middle chunk: flipping Bold on and back
    # off breaks the run at $b1/$b2 without altering the visible/resolved
    # formatting (keeps w:color=000000 on all three resulting runs).
    $middle = $d.Range($b1, $b2)
    $middle.Bold = 1
    $middle.Bold = 0

    Write-Output "split ok"
} else {
    Write-Output "target text not found"
}
